$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the population name text used in B2 (and shared by B4:B7) with flipped word order
$ws.Range("B2").Value = "ICER - ICER RRMM 2022 report - 12/19/2022"

# Update the population name text used in B3 with flipped word order
$ws.Range("B3").Value = "PRODFix_QOL_ECON - UtilityOutcome - 9/19/2022"

# B4:B7 reuse the same updated text as B2
$ws.Range("B4").Value = "ICER - ICER RRMM 2022 report - 12/19/2022"
$ws.Range("B5").Value = "ICER - ICER RRMM 2022 report - 12/19/2022"
$ws.Range("B6").Value = "ICER - ICER RRMM 2022 report - 12/19/2022"
$ws.Range("B7").Value = "ICER - ICER RRMM 2022 report - 12/19/2022"

# Update the active selection to B3 (was D3)
$ws.Range("B3").Select()
